$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A to text so the date string isn't reinterpreted as a date serial.
$ws.Range("A14").NumberFormat = "@"

# Append the newest mod-count data point as row 14.
$ws.Range("A14").Value = "2025/11/23"
$ws.Range("B14").Value = "逃离鸭科夫"
$ws.Range("C14").Value = 1236

# Mirror the formatting (centered style) used by the preceding data rows.
$ws.Range("A13:C13").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122)
